# Implement mixed languages on activity executions:
# add a new "Gemischte Durchfuehrung" (mixed execution) column (H) to the
# upload template, with sample values on the two sample data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in H1
$ws.Range("H1").Value = "Gemischte Durchfuehrung"

# Sample data in H2/H3, matching the existing "ja" answers in column G
$ws.Range("H2").Value = "ja"
$ws.Range("H3").Value = "ja"

# Row 3 (and the blank row 4) carry a distinct style (border/font) applied
# to the rest of that row - copy it onto the new column so H lines up.
$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H3").Value = "ja"

$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122) # xlPasteFormats

# Leave the selection on the newly added cell, like the saved workbook.
$ws.Range("H4").Select()
